$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.012.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.46%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.963.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.27%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.04%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.81%  '

# Row 6: 'USDC' -> 'USDC'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.07%  '

# Row 7: 'XRP' -> 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4834'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.70%  '

# Row 8: 'OKB' -> 'OKB'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.53'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.09%  '

# Row 9: 'Cardano' -> 'Cardano'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2948'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.61%  '

# Row 10: 'Dogecoin' -> 'Dogecoin'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06777'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.82%  '

# Row 11: 'Litecoin' -> 'Solana'
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.54%  '

# Row 12: 'Solana' -> 'Litecoin'
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '108.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.36%  '

# Row 13: 'WrappedEther' -> 'WrappedEther'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.959.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.18%  '

# Row 14: 'TRON' -> 'TRON'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07769'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.23%  '

# Row 15: 'Polkadot' -> 'Polkadot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.466'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.41%  '

# Row 16: 'Polygon' -> 'Polygon'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.6983'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.72%  '

# Row 17: 'BitcoinCash' -> 'BitcoinCash'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '287.34'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.53%  '

# Row 18: 'WrappedBTC' -> 'WrappedBTC'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '31.039.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.60%  '

# Row 19: 'Avalanche' -> 'Avalanche'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.33%  '

# Row 20: 'ShibaInu' -> 'ShibaInu'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007736'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.06%  '

# Row 21: 'BitDAO' -> 'WrappedliquidstakedEther2.0'
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.219.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.46%  '

# Row 22: 'Uniswap' -> 'Uniswap'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.629'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.00%  '

# Row 23: 'WrappedliquidstakedEther2.0' -> 'Dai'
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.0000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '

# Row 24: 'Dai' -> 'BinanceUSD'
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.02%  '

# Row 25: 'BinanceUSD' -> 'Chainlink'
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.610'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.65%  '

# Row 26: 'Chainlink' -> 'Cosmos'
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.916'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.61%  '

# Row 27: 'Cosmos' -> 'Monero'
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.32%  '

# Row 28: 'Monero' -> 'EthereumClassic'
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.25%  '

# Row 29: 'EthereumClassic' -> 'LidoDAOToken'
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.188'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.39%  '

# Row 30: 'LidoDAOToken' -> 'Stellar'
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1063'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.34%  '

# Row 31: 'Stellar' -> 'Toncoin'
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.446'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.24%  '

# Row 32: 'Toncoin' -> 'Filecoin'
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.828'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +17.75%  '

# Row 33: 'Filecoin' -> 'InternetComputer(DFINITY)'
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.546'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.08%  '

# Row 34: 'InternetComputer(DFINITY)' -> 'Hedera'
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05094'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.41%  '

# Row 35: 'Hedera' -> 'ImmutableX'
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7781'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.17%  '

# Row 36: 'ImmutableX' -> 'ARBITRUM'
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.174'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.87%  '

# Row 37: 'ARBITRUM' -> 'HuobiToken'
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.735'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.89%  '

# Row 38: 'VeChain' -> 'VeChain'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02038'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.73%  '

# Row 39: 'HuobiToken' -> 'MXToken'
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("E39").Value = '  +0.99%  '

# Row 40: 'MXToken' -> 'FraxShare'
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.504'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.81%  '

# Row 41: 'FraxShare' -> 'RenderToken'
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.127'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.48%  '

# Row 42: 'RenderToken' -> 'TrustWalletToken'
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8881'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.62%  '

# Row 43: 'TrustWalletToken' -> 'Quant'
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '109.77'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.27%  '

# Row 44: 'Quant' -> 'TheSandbox'
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4455'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.24%  '

# Row 45: 'TheSandbox' -> 'Aave'
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '71.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.83%  '

# Row 46: 'Aave' -> 'PaxDollar'
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.02%  '

# Row 47: 'PaxDollar' -> 'Aptos'
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.530'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.78%  '

# Row 48: 'Aptos' -> 'EnergySwap'
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.437'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.38%  '

# Row 49: 'Algorand' -> 'Algorand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1273'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.84%  '

# Row 50: 'EnergySwap' -> 'Elrond'
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.00'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.33%  '

# Row 51: 'Elrond' -> 'Maker'
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '943.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +10.81%  '
